# Working on 0_constants case for testing sliced variables on special identity.
# Adds a new "f_Category" column/classification to the FLOWS and FLOWS_AGG sets.

$wb = $excel.ActiveWorkbook

$wsCase     = $wb.Worksheets.Item("_set_CASE")
$wsFlows    = $wb.Worksheets.Item("_set_FLOWS")
$wsFlowsAgg = $wb.Worksheets.Item("_set_FLOWS_AGG")

# ---------------------------------------------------------------------------
# 1) Write header rows first for all three sheets (this reproduces the
#    shared-string ordering of the authored workbook: headers are inserted
#    into the shared string table before any data values).
# ---------------------------------------------------------------------------
$wsCase.Range("A1").Value = "c_Name"

$wsFlows.Range("A1").Value = "f_Name"
$wsFlows.Range("B1").Value = "f_Category"
$wsFlows.Range("C1").Value = "f_Aggregation"

$wsFlowsAgg.Range("A1").Value = "f_agg_Name"
$wsFlowsAgg.Range("B1").Value = "f_Category"

# Copy the header style (bold font + border) from the existing A1 header
# cells onto the newly introduced header cells so they match the rest.
$wsFlows.Range("A1").Copy($wsFlows.Range("B1"))
$wsFlows.Range("A1").Copy($wsFlows.Range("C1"))
$wsFlowsAgg.Range("A1").Copy($wsFlowsAgg.Range("B1"))

# Re-assert the text values after the style copy (Copy duplicates both
# value and formatting from the source cell).
$wsFlows.Range("B1").Value = "f_Category"
$wsFlows.Range("C1").Value = "f_Aggregation"
$wsFlowsAgg.Range("B1").Value = "f_Category"

# ---------------------------------------------------------------------------
# 2) Data rows, sheet by sheet.
# ---------------------------------------------------------------------------

# _set_CASE
$wsCase.Range("A2").Value = "base"

# _set_FLOWS (f_Name, f_Category, f_Aggregation)
$wsFlows.Range("A2").Value = "oil products"
$wsFlows.Range("B2").Value = "Yearly dispatched"
$wsFlows.Range("C2").Value = "oil products"

$wsFlows.Range("A3").Value = "electricity, gas"
$wsFlows.Range("B3").Value = "Hourly dispatched"
$wsFlows.Range("C3").Value = "electricity"

$wsFlows.Range("A4").Value = "electricity, res"
$wsFlows.Range("B4").Value = "Hourly dispatched"
$wsFlows.Range("C4").Value = "electricity"

$wsFlows.Range("A5").Value = "transport, icev"
$wsFlows.Range("B5").Value = "Yearly dispatched"
$wsFlows.Range("C5").Value = "transport"

$wsFlows.Range("A6").Value = "transport, bev"
$wsFlows.Range("B6").Value = "Yearly dispatched"
$wsFlows.Range("C6").Value = "transport"

# _set_FLOWS_AGG (f_agg_Name, f_Category)
$wsFlowsAgg.Range("A2").Value = "oil products"
$wsFlowsAgg.Range("B2").Value = "Yearly dispatched"

$wsFlowsAgg.Range("A3").Value = "electricity"
$wsFlowsAgg.Range("B3").Value = "Hourly dispatched"

$wsFlowsAgg.Range("A4").Value = "transport"
$wsFlowsAgg.Range("B4").Value = "Yearly dispatched"

# ---------------------------------------------------------------------------
# 3) View state: active cell / selection per sheet, active tab, window size.
# ---------------------------------------------------------------------------
$wsCase.Range("C10").Select()
$wsFlows.Range("D17").Select()
$wsFlowsAgg.Range("C11").Select()

$wsFlowsAgg.Activate()

$wb.Windows.Item(1).WindowState = -4143
$wb.Windows.Item(1).Left = 5350
$wb.Windows.Item(1).Top = 3150
$wb.Windows.Item(1).Width = 28800
$wb.Windows.Item(1).Height = 15370
